$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove workbook (structure) protection -> drops the empty <workbookProtection/> element
$wb.Unprotect()

# Match the window geometry recorded in the saved workbookView
$win = $excel.ActiveWindow
$win.Left = 240
$win.Top = 15
$win.Width = 16095
$win.Height = 9660

# Make sure the data sheet is the selected / active tab
$ws.Select()

# Restore Excel's default page margins (in inches)
$ws.PageSetup.LeftMargin = 0.7 * 72
$ws.PageSetup.RightMargin = 0.7 * 72
$ws.PageSetup.TopMargin = 0.75 * 72
$ws.PageSetup.BottomMargin = 0.75 * 72
$ws.PageSetup.HeaderMargin = 0.3 * 72
$ws.PageSetup.FooterMargin = 0.3 * 72
